$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column AM (39th column), shifting
# "MatchSequence" from AM1 to AP1, and fill the new columns with
# headers V39, V40, V41.
$ws.Range("AM1:AO1").EntireColumn.Insert()

$ws.Range("AM1").Value = "V39"
$ws.Range("AN1").Value = "V40"
$ws.Range("AO1").Value = "V41"
